# Applies the source-controlled update to "Artfynd" sheet rows 4-6.
# The three observation records (rows 4, 5, 6) were cyclically rotated:
#   new row4 = old row6, new row5 = old row4, new row6 = old row5
# (Id/TaxonId/species/etc. all move together; a few fields happen to keep
# the same value after rotation and are therefore left untouched below.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 <- previous Row 6
$ws.Range("A4").Value = 111782565
$ws.Range("B4").Value = 100532
$ws.Range("D4").Value = "CR"
$ws.Range("E4").Value = 223246
$ws.Range("F4").Value = "Skogsalm"
$ws.Range("G4").Value = "Ulmus glabra"
$ws.Range("H4").Value = "Huds."
$ws.Range("J4").Value = "plantor/tuvor"
$ws.Range("Q4").Value = 573877.2060252801
$ws.Range("R4").Value = 6303225.547499124

# Row 5 <- previous Row 4
$ws.Range("A5").Value = 111782566
$ws.Range("B5").Value = 103369
$ws.Range("E5").Value = 221423
$ws.Range("F5").Value = "Myskmadra"
$ws.Range("G5").Value = "Galium odoratum"
$ws.Range("H5").Value = "(L.) Scop."
$ws.Range("I5").Value = "1"
$ws.Range("J5").Value = "m²"
$ws.Range("Q5").Value = 573877.0511306904
$ws.Range("R5").Value = 6303234.29156004

# Row 6 <- previous Row 5
$ws.Range("A6").Value = 111782567
$ws.Range("B6").Value = 98535
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 222498
$ws.Range("F6").Value = "Blåsippa"
$ws.Range("G6").Value = "Hepatica nobilis"
$ws.Range("H6").Value = "Schreb."
$ws.Range("I6").Value = "5"
$ws.Range("Q6").Value = 573909.350056502
$ws.Range("R6").Value = 6303235.410511858
